$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1728506666666667
$ws.Range("H2").Value = 0.518552
$ws.Range("I2").Value = 0.0840503369699626
$ws.Range("J2").Value = 0.0840503369699626
$ws.Range("M2").Value = 19.163974
$ws.Range("N2").Value = 57.491922
$ws.Range("O2").Value = 0.6845732287637933
$ws.Range("P2").Value = 0.6845732287637933
$ws.Range("Q2").Value = 3.312505681882667
$ws.Range("R2").Value = 29.812551136944
$ws.Range("S2").Value = 0.05753861055821213
$ws.Range("T2").Value = 0.05753861055821213

$ws.Range("G3").Value = 0.1728506666666667
$ws.Range("H3").Value = 0.518552
$ws.Range("I3").Value = 0.0840503369699626
$ws.Range("J3").Value = 0.0840503369699626
$ws.Range("O3").Value = 0.02733363438148322
$ws.Range("P3").Value = 0.02733363438148323
$ws.Range("Q3").Value = 0.1322617002693333
$ws.Range("R3").Value = 1.190355302424
$ws.Range("S3").Value = 0.00229740118037742
$ws.Range("T3").Value = 0.002297401180377421

$ws.Range("G4").Value = 0.1728506666666667
$ws.Range("H4").Value = 0.518552
$ws.Range("I4").Value = 0.0840503369699626
$ws.Range("J4").Value = 0.0840503369699626
$ws.Range("M4").Value = 7.880893333333333
$ws.Range("N4").Value = 23.64268
$ws.Range("O4").Value = 0.281520346184098
$ws.Range("P4").Value = 0.281520346184098
$ws.Range("Q4").Value = 1.362217666595556
$ws.Range("R4").Value = 12.25995899936
$ws.Range("S4").Value = 0.02366187996067396
$ws.Range("T4").Value = 0.02366187996067396

$ws.Range("G5").Value = 0.1728506666666667
$ws.Range("H5").Value = 0.518552
$ws.Range("I5").Value = 0.0840503369699626
$ws.Range("J5").Value = 0.0840503369699626
$ws.Range("M5").Value = 0.183999
$ws.Range("N5").Value = 0.551997
$ws.Range("O5").Value = 0.006572790670625477
$ws.Range("P5").Value = 0.006572790670625476
$ws.Range("Q5").Value = 0.031804349816
$ws.Range("R5").Value = 0.286239148344
$ws.Range("S5").Value = 0.0005524452706990978
$ws.Range("T5").Value = 0.0005524452706990977

$ws.Range("I6").Value = 0.6650661694281633
$ws.Range("J6").Value = 0.6650661694281633
$ws.Range("M6").Value = 19.163974
$ws.Range("N6").Value = 57.491922
$ws.Range("O6").Value = 0.6845732287637933
$ws.Range("P6").Value = 0.6845732287637933
$ws.Range("Q6").Value = 26.21090580334067
$ws.Range("R6").Value = 235.898152230066
$ws.Range("S6").Value = 0.4552864949470058
$ws.Range("T6").Value = 0.4552864949470058

$ws.Range("I7").Value = 0.6650661694281633
$ws.Range("J7").Value = 0.6650661694281633
$ws.Range("O7").Value = 0.02733363438148322
$ws.Range("P7").Value = 0.02733363438148323
$ws.Range("S7").Value = 0.01817867551464299
$ws.Range("T7").Value = 0.01817867551464299

$ws.Range("I8").Value = 0.6650661694281633
$ws.Range("J8").Value = 0.6650661694281633
$ws.Range("M8").Value = 7.880893333333333
$ws.Range("N8").Value = 23.64268
$ws.Range("O8").Value = 0.281520346184098
$ws.Range("P8").Value = 0.281520346184098
$ws.Range("Q8").Value = 10.77883704111555
$ws.Range("R8").Value = 97.00953337003999
$ws.Range("S8").Value = 0.1872296582527485
$ws.Range("T8").Value = 0.1872296582527485

$ws.Range("I9").Value = 0.6650661694281633
$ws.Range("J9").Value = 0.6650661694281633
$ws.Range("M9").Value = 0.183999
$ws.Range("N9").Value = 0.551997
$ws.Range("O9").Value = 0.006572790670625477
$ws.Range("P9").Value = 0.006572790670625476
$ws.Range("Q9").Value = 0.251658682949
$ws.Range("R9").Value = 2.264928146541
$ws.Range("S9").Value = 0.004371340713766055
$ws.Range("T9").Value = 0.004371340713766055

$ws.Range("G10").Value = 0.5159453333333334
$ws.Range("H10").Value = 1.547836
$ws.Range("I10").Value = 0.2508834936018741
$ws.Range("J10").Value = 0.2508834936018741
$ws.Range("M10").Value = 19.163974
$ws.Range("N10").Value = 57.491922
$ws.Range("O10").Value = 0.6845732287637933
$ws.Range("P10").Value = 0.6845732287637933
$ws.Range("Q10").Value = 9.887562953421334
$ws.Range("R10").Value = 88.98806658079201
$ws.Range("S10").Value = 0.1717481232585754
$ws.Range("T10").Value = 0.1717481232585754

$ws.Range("G11").Value = 0.5159453333333334
$ws.Range("H11").Value = 1.547836
$ws.Range("I11").Value = 0.2508834936018741
$ws.Range("J11").Value = 0.2508834936018741
$ws.Range("O11").Value = 0.02733363438148322
$ws.Range("P11").Value = 0.02733363438148323
$ws.Range("Q11").Value = 0.3947905342146666
$ws.Range("R11").Value = 3.553114807932
$ws.Range("S11").Value = 0.006857557686462812
$ws.Range("T11").Value = 0.006857557686462813

$ws.Range("G12").Value = 0.5159453333333334
$ws.Range("H12").Value = 1.547836
$ws.Range("I12").Value = 0.2508834936018741
$ws.Range("J12").Value = 0.2508834936018741
$ws.Range("M12").Value = 7.880893333333333
$ws.Range("N12").Value = 23.64268
$ws.Range("O12").Value = 0.281520346184098
$ws.Range("P12").Value = 0.281520346184098
$ws.Range("Q12").Value = 4.066110137831111
$ws.Range("R12").Value = 36.59499124048
$ws.Range("S12").Value = 0.07062880797067554
$ws.Range("T12").Value = 0.07062880797067554

$ws.Range("G13").Value = 0.5159453333333334
$ws.Range("H13").Value = 1.547836
$ws.Range("I13").Value = 0.2508834936018741
$ws.Range("J13").Value = 0.2508834936018741
$ws.Range("M13").Value = 0.183999
$ws.Range("N13").Value = 0.551997
$ws.Range("O13").Value = 0.006572790670625477
$ws.Range("P13").Value = 0.006572790670625476
$ws.Range("Q13").Value = 0.09493342538800001
$ws.Range("R13").Value = 0.8544008284919999
$ws.Range("S13").Value = 0.001649004686160325
$ws.Range("T13").Value = 0.001649004686160325
